$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I3").Value = "sv"
$ws.Range("J3").Value = "Statement-opinion"

# Row 11: aa/Agree/Accept -> sd/Statement-non-opinion
$ws.Range("I11").Value = "sd"
$ws.Range("J11").Value = "Statement-non-opinion"

# Row 16: sv/Statement-opinion -> sd/Statement-non-opinion
$ws.Range("I16").Value = "sd"
$ws.Range("J16").Value = "Statement-non-opinion"

# Row 21: sd/Statement-non-opinion -> sv/Statement-opinion
$ws.Range("I21").Value = "sv"
$ws.Range("J21").Value = "Statement-opinion"

# Row 43: qy/Yes-No-Question -> sv/Statement-opinion
$ws.Range("I43").Value = "sv"
$ws.Range("J43").Value = "Statement-opinion"

# Row 46: ba/Appreciation -> sd/Statement-non-opinion
$ws.Range("I46").Value = "sd"
$ws.Range("J46").Value = "Statement-non-opinion"
